$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New TPM-derived values for the Gal -> Gpr151 LR-pair table (rows 2-17).
# Each entry maps a cell reference to its updated value, taken from the
# recomputed NATMI statistics after the upstream TPM matrix was refreshed.
$updates = [ordered]@{
    "E2" = 3
    "F2" = 1
    "G2" = 1.776285
    "H2" = 5.328855000000001
    "I2" = 0.2250740306326953
    "J2" = 0.2250740306326953
    "M2" = 1.059024
    "N2" = 3.177072
    "O2" = 0.401751683632714
    "P2" = 0.401751683632714
    "Q2" = 1.88112844584
    "R2" = 16.93015601256
    "S2" = 0.09042387074868638
    "T2" = 0.09042387074868638
    "E3" = 3
    "F3" = 1
    "G3" = 1.776285
    "H3" = 5.328855000000001
    "I3" = 0.2250740306326953
    "J3" = 0.2250740306326953
    "O3" = 0.1572783628427189
    "P3" = 0.1572783628427188
    "Q3" = 0.7364270376750002
    "R3" = 6.627843339075001
    "S3" = 0.03539927505632227
    "T3" = 0.03539927505632227
    "E4" = 3
    "F4" = 1
    "G4" = 1.776285
    "H4" = 5.328855000000001
    "I4" = 0.2250740306326953
    "J4" = 0.2250740306326953
    "M4" = 0.9285883333333333
    "N4" = 2.785765
    "O4" = 0.3522695673736974
    "P4" = 0.3522695673736974
    "Q4" = 1.649437527675
    "R4" = 14.844937749075
    "S4" = 0.0792867313980339
    "T4" = 0.07928673139803388
    "E5" = 3
    "F5" = 1
    "G5" = 1.776285
    "H5" = 5.328855000000001
    "I5" = 0.2250740306326953
    "J5" = 0.2250740306326953
    "M5" = 0.2338156666666666
    "N5" = 0.7014469999999999
    "O5" = 0.08870038615086982
    "P5" = 0.08870038615086982
    "Q5" = 0.415323261465
    "R5" = 3.737909353185
    "S5" = 0.01996415342965278
    "T5" = 0.01996415342965278
    "I6" = 0.1783607964348723
    "J6" = 0.1783607964348722
    "M6" = 1.059024
    "N6" = 3.177072
    "O6" = 0.401751683632714
    "P6" = 0.401751683632714
    "Q6" = 1.490707598976
    "R6" = 13.416368390784
    "S6" = 0.0716567502617817
    "T6" = 0.07165675026178168
    "I7" = 0.1783607964348723
    "J7" = 0.1783607964348722
    "O7" = 0.1572783628427189
    "P7" = 0.1572783628427188
    "S7" = 0.02805229405860016
    "T7" = 0.02805229405860015
    "I8" = 0.1783607964348723
    "J8" = 0.1783607964348722
    "M8" = 0.9285883333333333
    "N8" = 2.785765
    "O8" = 0.3522695673736974
    "P8" = 0.3522695673736974
    "Q8" = 1.30710322412
    "R8" = 11.76392901708
    "S8" = 0.06283108059654056
    "T8" = 0.06283108059654055
    "I9" = 0.1783607964348723
    "J9" = 0.1783607964348722
    "M9" = 0.2338156666666666
    "N9" = 0.7014469999999999
    "O9" = 0.08870038615086982
    "P9" = 0.08870038615086982
    "Q9" = 0.3291245439759999
    "R9" = 2.962120895784
    "S9" = 0.01582067151794985
    "T9" = 0.01582067151794985
    "G10" = 4.248107333333333
    "H10" = 12.744322
    "I10" = 0.5382799720054182
    "J10" = 0.5382799720054181
    "M10" = 1.059024
    "N10" = 3.177072
    "O10" = 0.401751683632714
    "P10" = 0.401751683632714
    "Q10" = 4.498847620576
    "R10" = 40.489628585184
    "S10" = 0.2162548850189469
    "T10" = 0.2162548850189469
    "G11" = 4.248107333333333
    "H11" = 12.744322
    "I11" = 0.5382799720054182
    "J11" = 0.5382799720054181
    "O11" = 0.1572783628427189
    "P11" = 0.1572783628427188
    "Q11" = 1.761215739147778
    "R11" = 15.85094165233
    "S11" = 0.08465979274803673
    "T11" = 0.08465979274803669
    "G12" = 4.248107333333333
    "H12" = 12.744322
    "I12" = 0.5382799720054182
    "J12" = 0.5382799720054181
    "M12" = 0.9285883333333333
    "N12" = 2.785765
    "O12" = 0.3522695673736974
    "P12" = 0.3522695673736974
    "Q12" = 3.944742908481111
    "R12" = 35.50268617633
    "S12" = 0.1896196528642746
    "T12" = 0.1896196528642746
    "G13" = 4.248107333333333
    "H13" = 12.744322
    "I13" = 0.5382799720054182
    "J13" = 0.5382799720054181
    "M13" = 0.2338156666666666
    "N13" = 0.7014469999999999
    "O13" = 0.08870038615086982
    "P13" = 0.08870038615086982
    "Q13" = 0.9932740482148888
    "R13" = 8.939466433933999
    "S13" = 0.04774564137415999
    "T13" = 0.04774564137415999
    "E14" = 2
    "F14" = 0.6666666666666666
    "G14" = 0.459987
    "H14" = 1.379961
    "I14" = 0.05828520092701431
    "J14" = 0.0582852009270143
    "M14" = 1.059024
    "N14" = 3.177072
    "O14" = 0.401751683632714
    "P14" = 0.401751683632714
    "Q14" = 0.4871372726879999
    "R14" = 4.384235454192
    "S14" = 0.02341617760329902
    "T14" = 0.02341617760329901
    "E15" = 2
    "F15" = 0.6666666666666666
    "G15" = 0.459987
    "H15" = 1.379961
    "I15" = 0.05828520092701431
    "J15" = 0.0582852009270143
    "O15" = 0.1572783628427189
    "P15" = 0.1572783628427188
    "Q15" = 0.190705243685
    "R15" = 1.716347193165
    "S15" = 0.009167000979759731
    "T15" = 0.009167000979759728
    "E16" = 2
    "F16" = 0.6666666666666666
    "G16" = 0.459987
    "H16" = 1.379961
    "I16" = 0.05828520092701431
    "J16" = 0.0582852009270143
    "M16" = 0.9285883333333333
    "N16" = 2.785765
    "O16" = 0.3522695673736974
    "P16" = 0.3522695673736974
    "Q16" = 0.427138561685
    "R16" = 3.844247055165
    "S16" = 0.02053210251484836
    "T16" = 0.02053210251484835
    "E17" = 2
    "F17" = 0.6666666666666666
    "G17" = 0.459987
    "H17" = 1.379961
    "I17" = 0.05828520092701431
    "J17" = 0.0582852009270143
    "M17" = 0.2338156666666666
    "N17" = 0.7014469999999999
    "O17" = 0.08870038615086982
    "P17" = 0.08870038615086982
    "Q17" = 0.107552167063
    "R17" = 0.967969503567
    "S17" = 0.005169919829107205
    "T17" = 0.005169919829107204
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
